$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (uppercased station / stop names, expanded service type) ---
$ws.Range("B2").Value  = "VIDYA NAGAR"
$ws.Range("B5").Value  = "VIDYA NAGAR,HUBLI"
$ws.Range("B7").Value  = "VIDYA NAGAR 1"
$ws.Range("B9").Value  = "City,Suburban,All Stops"

# --- New column headers C/D ---
$ws.Range("C1").Value = "C"
$ws.Range("D1").Value = "D"

# --- New credential / url block in columns C/D ---
$ws.Range("C2").Value = "URL"
$ws.Range("D2").Value = "http://nechubli.com:5001"

$ws.Range("C3").Value = "User_Name"
$ws.Range("D3").Value = "VPSD"

$ws.Range("C4").Value = "Password"
$ws.Range("D4").Value = "Vpsd@master@123"

# --- New rows below the existing table ---
$ws.Range("A12").Value = "Color Code"
$ws.Range("B12").Value = "Green"

$ws.Range("A14").Value = "Station_Name"
$ws.Range("B14").Value = "VIDYA NAGAR BVB"

# --- Formatting: center-align the whole used area ---
$ws.Range("A1:B12").HorizontalAlignment = -4108
$ws.Range("C1:D5").HorizontalAlignment = -4108
$ws.Range("A14").HorizontalAlignment = -4108

# --- Borders around the main A:B table (rows 1-12) ---
$ws.Range("A1:B12").Borders.LineStyle = 1

# --- Hyperlinks for the URL / Password cells (added after alignment so the
#     built-in Hyperlink style inherits the already-centered xf) ---
$ws.Hyperlinks.Add($ws.Range("D2"), "http://nechubli.com:5001")
$ws.Hyperlinks.Add($ws.Range("D4"), "http://nechubli.com:5001")

# --- Column widths (approximate best-fit sizing) ---
$ws.Columns(1).ColumnWidth = 18.6
$ws.Columns(2).ColumnWidth = 20.6
$ws.Columns(4).ColumnWidth = 23.1

# --- Selection mirrors the author's final cursor position ---
$ws.Range("A15").Select() | Out-Null

Write-Host "done"
